# Implemented getting kafka relations.
#
# 1. Add a "Weight" header to the existing feignRelations sheet (column G,
#    row 1) — the G column values themselves already existed.
# 2. Add a new "kafkaRelations" worksheet with its own header row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("feignRelations")

# 1. Add missing header to feignRelations!G1
$ws1.Range("G1").Value = "Weight"

# 2. Create the new kafkaRelations sheet (inserted after feignRelations)
#    and populate its header row.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "kafkaRelations"

$ws2.Range("A1").Value = "Source Service Name"
$ws2.Range("B1").Value = "Source Class Name"
$ws2.Range("C1").Value = "Source Method Signature"
$ws2.Range("D1").Value = "Target Service Name"
$ws2.Range("E1").Value = "Target Class Name"
$ws2.Range("F1").Value = "Target Method Signature"
$ws2.Range("G1").Value = "Topic"
$ws2.Range("H1").Value = "Weight"
